$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 20
$ws.Range("B1").Value = 23
$ws.Range("C1").Value = 25
$ws.Range("D1").Value = 34
$ws.Range("E1").Value = 39
$ws.Range("F1").Value = 44
$ws.Range("G1").Value = "과거기록 : [194]회차 4등"

# Row 2
$ws.Range("A2").Value = 13
$ws.Range("B2").Value = 18
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 21
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = 36
$ws.Range("G2").Value = "과거기록 : [246, 64]회차 4등"

# Row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 18
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 34
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 45
$ws.Range("G3").Value = "과거기록 : [616, 396, 183, 94]회차 4등"

# Row 4
$ws.Range("A4").Value = 19
$ws.Range("B4").Value = 23
$ws.Range("C4").Value = 26
$ws.Range("D4").Value = 27
$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 33
$ws.Range("G4").Value = "과거기록 : [955, 570, 530, 433]회차 4등"

# Row 5
$ws.Range("A5").Value = 14
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 34
$ws.Range("F5").Value = 41
$ws.Range("G5").Value = "과거기록 : [306]회차 4등"
